# Fix clean_data_central_africa.xlsx: reset index before saving.
# This:
#   1. Resets the country_index column (A) to a sequential 0-based index
#      (the old "global" index values are replaced by the row's position).
#   2. Also fixes a data value (compulsory_edu_yrs for row 2) that changed
#      as part of recomputing the export.
#   3. Reorders the gni_index column so it comes after
#      pct_industry_employment / pct_services_employment / exports_pct_gdp /
#      fdi_pct_gdp (columns O:S), matching the corrected column order.
#   4. Drops four stray "m_" columns that shouldn't have been exported
#      (m_income_group, m_homicides_per_100k, m_adult_literacy_pct,
#      m_tax_revenue_pct_gdp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Reset the country_index column (A2:A15) to 0..13 -------------------
$indexRows = 14
$indexVals = New-Object 'object[,]' $indexRows,1
for ($i = 0; $i -lt $indexRows; $i++) {
    $indexVals[$i,0] = $i
}
$ws.Range("A2:A15").Value = $indexVals

# --- 2. Data correction: compulsory_edu_yrs for Burundi (row 2) ------------
$ws.Range("K2").Value = 9.755208333333334

# --- 3. Reorder O:S (gni_index moves after the other four columns) ---------
$reorderRows = 14
$oldVals = $ws.Range("O2:S15").Value()
$newVals = New-Object 'object[,]' $reorderRows,5
for ($i = 1; $i -le $reorderRows; $i++) {
    $gniIndex   = $oldVals[$i,1]
    $industry   = $oldVals[$i,2]
    $services   = $oldVals[$i,3]
    $exports    = $oldVals[$i,4]
    $fdi        = $oldVals[$i,5]
    $newVals[$i-1,0] = $industry
    $newVals[$i-1,1] = $services
    $newVals[$i-1,2] = $exports
    $newVals[$i-1,3] = $fdi
    $newVals[$i-1,4] = $gniIndex
}
$ws.Range("O2:S15").Value = $newVals

# Matching header reorder for row 1.
$ws.Range("O1").Value = "pct_industry_employment"
$ws.Range("P1").Value = "pct_services_employment"
$ws.Range("Q1").Value = "exports_pct_gdp"
$ws.Range("R1").Value = "fdi_pct_gdp"
$ws.Range("S1").Value = "gni_index"

# --- 4. Drop the four stray "m_" columns ------------------------------------
# Delete from rightmost to leftmost so earlier deletions don't shift the
# column letters of the ones still to be removed.
$ws.Columns("AZ").Delete()   # m_tax_revenue_pct_gdp
$ws.Columns("AV").Delete()   # m_adult_literacy_pct
$ws.Columns("AU").Delete()   # m_homicides_per_100k
$ws.Columns("AD").Delete()   # m_income_group
